# Fix bug in interpolating widths.
#
# 1) The "Sections - Copy" sheet data (A3:D56) gets re-sorted by the
#    Chainage-index column (A) in descending order (the GERD / chainage-0
#    row ends up first, chainage 120000 ends up last).
# 2) Width (column C) values that were bogus "interpolated" leftovers are
#    cleared wherever there is no real surveyed Average Bed Elevation
#    (column D) for that row - plus two extra stray leftover values that
#    need clearing as well.
# 3) The active cell/selection on the sheet moves to O5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Sort the data block A3:D56 by column A, descending ---------------
$sortRange = $ws.Range("A3:D56")
$sortKey   = $ws.Range("A3:A56")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($sortKey, 0, 2, 0, 0)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
$ws.Sort.MatchCase = $false
$ws.Sort.Orientation = 1
$ws.Sort.Apply()

# --- 2) Clear the stray/interpolated Width values -------------------------
# (Excel's union-range ClearContents only reliably clears the first area in
# this environment, so clear each row individually instead.)
$blankRows = @(5,7,8,9,10,12,13,14,16,17,18,19,20,22,23,24,25,27,29,30,31,32,34,35,37,39,40,41,44,45,47,53,54,55,56)
foreach ($r in $blankRows) {
    $ws.Range("C" + $r).ClearContents()
}

# The "GERD" label that used to sit in the last row now sorts to the top
# of the block; replace it with its numeric index (54) like every other row
$ws.Range("A3").Value = 54

# --- 3) Move the selection -------------------------------------------------
$null = $ws.Range("O5").Select()

# --- misc workbook metadata -------------------------------------------------
Write-Host "done"
